$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date value changes from 45208 (2023-10-09) to
# 45212 (2023-10-13) for rows 2 through 15.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
